# Updates the cryptos price list (columns D = Price, E = Volume(1h))
# to the latest scraped values, mirroring the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.128.31"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "2.399.78"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'504.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.62%  "
$ws.Range("D6").Value = "'132.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.37%  "
$ws.Range("D7").Value = "'0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("D8").Value = "'0.559"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("D9").Value = "2.433.46"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").Value = "'0.0974"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "'4.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.63%  "
$ws.Range("D14").Value = "2.835.84"
$ws.Range("E14").Value = "  -1.56%  "
$ws.Range("D15").Value = "57.062.23"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").Value = "'21.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "'0.0000134"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").Value = "2.402.72"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("D19").Value = "'10.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").Value = "'4.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").Value = "'313.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D22").Value = "'6.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.71%  "
$ws.Range("D23").Value = "'0.995"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("E24").Value = "  -2.04%  "
$ws.Range("D25").Value = "'65.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").Value = "'0.994"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("D27").Value = "2.511.22"
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("D28").Value = "'0.382"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.40%  "
$ws.Range("E29").Value = "  -3.03%  "
$ws.Range("D30").Value = "'7.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.27%  "
$ws.Range("D31").Value = "'172.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").Value = "0.0₃0732"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").Value = "'1.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("D34").Value = "'6.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("D35").Value = "'1.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.99%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("D38").Value = "'18.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.01%  "
$ws.Range("E39").Value = "  +3.85%  "
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("D41").Value = "'0.815"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.36%  "
$ws.Range("D42").Value = "'36.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("D44").Value = "'132.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.82%  "
$ws.Range("D45").Value = "'3.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "'5.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.48%  "
$ws.Range("D47").Value = "'257.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.62%  "
$ws.Range("E48").Value = "  -2.45%  "
$ws.Range("D49").Value = "'0.0915"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("D50").Value = "'0.0495"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("D51").Value = "'0.0212"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.97%  "
